$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the "food sets count" (套餐销售组数) row (row 29) with placeholder
# template tokens, one per store column B..F.
$ws.Range("B29").Value = '${TS_foodSetsCount}'
$ws.Range("C29").Value = '${JN_foodSetsCount}'
$ws.Range("D29").Value = '${BY_foodSetsCount}'
$ws.Range("E29").Value = '${HX_foodSetsCount}'
$ws.Range("F29").Value = '${DY_foodSetsCount}'

# Update the view state to reflect where the user left the selection/scroll.
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Range("F30").Select()
